$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "blockchain free" / "block.chain.technology" -> "helix jump" / "com.singleton.helix"
$ws.Range("A4").Value2 = "helix jump"
$ws.Range("B4").Value2 = "com.singleton.helix"

# B4 needs the same (unwrapped / default) formatting that A4 / A5 / B5 already use,
# instead of the monospace wrapped style it had before. Copy formats only from A4.
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# Move the active selection to A4:B4
[void]$ws.Range("A4:B4").Select()
